$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data cells for rows 2-5 (daily satellite pass data) ---
# Row 2
$ws.Range("A2").Value = "20260221--01"
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = "04:00"
$ws.Range("D2").Value = "00:00"
$ws.Range("E2").Value = "05:07:07"
$ws.Range("F2").Value = "05:10:00"
$ws.Range("G2").Value = "05:12:00"
$ws.Range("H2").Value = "05:14:00"
$ws.Range("I2").Value = "05:16:53"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = -9.1
$ws.Range("M2").Value = "A+B"
$ws.Range("N2").Value = "2"
$ws.Range("O2").Value = 100
$ws.Range("P2").Value = 98
$ws.Range("Q2").Value = 94
$ws.Range("R2").Value = 68

# Row 3
$ws.Range("A3").Value = "20260222--01"
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = "03:07"
$ws.Range("D3").Value = "03:07"
$ws.Range("E3").Value = "04:28:46"
$ws.Range("F3").Value = "04:31:55"
$ws.Range("G3").Value = "04:33:28"
$ws.Range("H3").Value = "04:35:02"
$ws.Range("I3").Value = "04:38:11"
$ws.Range("J3").Value = "4°"
$ws.Range("K3").Value = "04:30:04"
$ws.Range("L3").Value = -15
$ws.Range("M3").Value = "A"
$ws.Range("N3").Value = "4"
$ws.Range("O3").Value = 100
$ws.Range("P3").Value = 99
$ws.Range("Q3").Value = 98
$ws.Range("R3").Value = 100

# Row 4
$ws.Range("A4").Value = "20260223--01"
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "01:30"
$ws.Range("D4").Value = "01:30"
$ws.Range("E4").Value = "03:50:24"
$ws.Range("F4").Value = "03:54:09"
$ws.Range("G4").Value = "03:54:54"
$ws.Range("H4").Value = "03:55:39"
$ws.Range("I4").Value = "03:59:25"
$ws.Range("J4").Value = "9°"
$ws.Range("K4").Value = "03:53:25"
$ws.Range("L4").Value = -20.8
$ws.Range("M4").Value = "A"
$ws.Range("N4").Value = "4"
$ws.Range("O4").Value = 100
$ws.Range("P4").Value = 57
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = 100

# Row 5
$ws.Range("A5").Value = "20260224--01"
$ws.Range("B5").Value = 24
$ws.Range("C5").Value = "05:51"
$ws.Range("D5").Value = "05:51"
$ws.Range("E5").Value = "04:48:07"
$ws.Range("F5").Value = "04:50:30"
$ws.Range("G5").Value = "04:53:25"
$ws.Range("H5").Value = "04:56:21"
$ws.Range("I5").Value = "04:58:45"
$ws.Range("J5").Value = "8°"
$ws.Range("K5").Value = "04:50:08"
$ws.Range("L5").Value = -11.2
$ws.Range("M5").Value = "A"
$ws.Range("N5").Value = "2"
$ws.Range("O5").Value = 88
$ws.Range("P5").Value = 64
$ws.Range("Q5").Value = 39
$ws.Range("R5").Value = 19

# --- Update heatmap fill/font colors for columns O:R (rows 2-5) ---
$c = $ws.Range("O2")
$c.Interior.Color = 8351984
$c.Font.Color = 3355443
$c = $ws.Range("P2")
$c.Interior.Color = 16777215
$c.Font.Color = 3355443
$c = $ws.Range("Q2")
$c.Interior.Color = 16579062
$c.Font.Color = 3355443
$c = $ws.Range("R2")
$c.Interior.Color = 15654092
$c.Font.Color = 3355443
$c = $ws.Range("O3")
$c.Interior.Color = 8351984
$c.Font.Color = 3355443
$c = $ws.Range("P3")
$c.Interior.Color = 16777215
$c.Font.Color = 3355443
$c = $ws.Range("Q3")
$c.Interior.Color = 16777215
$c.Font.Color = 3355443
$c = $ws.Range("R3")
$c.Interior.Color = 16777215
$c.Font.Color = 3355443
$c = $ws.Range("O4")
$c.Interior.Color = 8351984
$c.Font.Color = 3355443
$c = $ws.Range("P4")
$c.Interior.Color = 15125426
$c.Font.Color = 3355443
$c = $ws.Range("Q4")
$c.Interior.Color = 13671790
$c.Font.Color = 3355443
$c = $ws.Range("R4")
$c.Interior.Color = 16777215
$c.Font.Color = 3355443
$c = $ws.Range("O5")
$c.Interior.Color = 8351984
$c.Font.Color = 3355443
$c = $ws.Range("P5")
$c.Interior.Color = 15521732
$c.Font.Color = 3355443
$c = $ws.Range("Q5")
$c.Interior.Color = 14596761
$c.Font.Color = 3355443
$c = $ws.Range("R5")
$c.Interior.Color = 13869943
$c.Font.Color = 3355443

# --- Extend conditional formatting ranges from row 4 to row 5 ---
$rng = $ws.Range("A2:A4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("A2:A5"))
}
$rng = $ws.Range("B2:B4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("B2:B5"))
}
$rng = $ws.Range("C2:C4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("C2:C5"))
}
$rng = $ws.Range("D2:D4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("D2:D5"))
}
$rng = $ws.Range("E2:E4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("E2:E5"))
}
$rng = $ws.Range("F2:F4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("F2:F5"))
}
$rng = $ws.Range("G2:G4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("G2:G5"))
}
$rng = $ws.Range("H2:H4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("H2:H5"))
}
$rng = $ws.Range("I2:I4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("I2:I5"))
}
$rng = $ws.Range("J2:J4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("J2:J5"))
}
$rng = $ws.Range("K2:K4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("K2:K5"))
}
$rng = $ws.Range("L2:L4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("L2:L5"))
}
$rng = $ws.Range("M2:M4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("M2:M5"))
}
$rng = $ws.Range("N2:N4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("N2:N5"))
}
$rng = $ws.Range("O2:O4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("O2:O5"))
}
$rng = $ws.Range("P2:P4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("P2:P5"))
}
$rng = $ws.Range("Q2:Q4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("Q2:Q5"))
}
$rng = $ws.Range("R2:R4")
for ($i = 1; $i -le $rng.FormatConditions.Count; $i++) {
    $rng.FormatConditions.Item($i).ModifyAppliesToRange($ws.Range("R2:R5"))
}

Write-Output "Edit complete"
